# Apply the Dual Blades motion-value sheet update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,17
$data[0,0] = 'Name'
$data[0,1] = 'Attack 1'
$data[0,2] = 'Attack 2'
$data[0,3] = 'Attack 3'
$data[0,4] = 'Attack 4'
$data[0,5] = 'Attack 5'
$data[0,6] = 'Attack 6'
$data[0,7] = 'Attack 7'
$data[0,8] = 'Attack 8'
$data[0,9] = 'Attack 9'
$data[0,10] = 'Attack 10'
$data[0,11] = 'Attack 11'
$data[0,12] = 'Attack 12'
$data[0,13] = 'Attack 13'
$data[0,14] = 'Attack 14'
$data[0,15] = 'Attack 15'
$data[0,16] = 'Attack 16'
$data[1,0] = 'Double Slash'
$data[1,1] = 8
$data[1,2] = 10
$data[1,3] = 0
$data[1,4] = 0
$data[1,5] = 0
$data[1,6] = 0
$data[1,7] = 0
$data[1,8] = 0
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 0
$data[1,14] = 0
$data[1,15] = 0
$data[1,16] = 0
$data[2,0] = 'Double Slash Return'
$data[2,1] = 9
$data[2,2] = 10
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 0
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 0
$data[2,14] = 0
$data[2,15] = 0
$data[2,16] = 0
$data[3,0] = 'Circle Slash'
$data[3,1] = 11
$data[3,2] = 9
$data[3,3] = 11
$data[3,4] = 0
$data[3,5] = 0
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 0
$data[3,14] = 0
$data[3,15] = 0
$data[3,16] = 0
$data[4,0] = 'Lunging Strike'
$data[4,1] = 3
$data[4,2] = 3
$data[4,3] = 5
$data[4,4] = 5
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = 0
$data[4,8] = 0
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 0
$data[4,14] = 0
$data[4,15] = 0
$data[4,16] = 0
$data[5,0] = 'Left Round Slash'
$data[5,1] = 7
$data[5,2] = 5
$data[5,3] = 15
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 0
$data[5,8] = 0
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 0
$data[5,16] = 0
$data[6,0] = 'Right Round Slash'
$data[6,1] = 7
$data[6,2] = 5
$data[6,3] = 15
$data[6,4] = 0
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 0
$data[6,14] = 0
$data[6,15] = 0
$data[6,16] = 0
$data[7,0] = 'Turn Slash'
$data[7,1] = 5
$data[7,2] = 5
$data[7,3] = 0
$data[7,4] = 0
$data[7,5] = 0
$data[7,6] = 0
$data[7,7] = 0
$data[7,8] = 0
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 0
$data[7,14] = 0
$data[7,15] = 0
$data[7,16] = 0
$data[8,0] = 'Rising Slash'
$data[8,1] = 7
$data[8,2] = 0
$data[8,3] = 0
$data[8,4] = 0
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 0
$data[8,8] = 0
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 0
$data[8,14] = 0
$data[8,15] = 0
$data[8,16] = 0
$data[9,0] = 'Sliding Slash'
$data[9,1] = 10
$data[9,2] = 10
$data[9,3] = 13
$data[9,4] = 13
$data[9,5] = 0
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 0
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 0
$data[9,16] = 0
$data[10,0] = 'Demon Fangs'
$data[10,1] = 11
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 0
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 0
$data[10,16] = 0
$data[11,0] = 'Twofold Demon Slash'
$data[11,1] = 7
$data[11,2] = 14
$data[11,3] = 0
$data[11,4] = 0
$data[11,5] = 0
$data[11,6] = 0
$data[11,7] = 0
$data[11,8] = 0
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 0
$data[11,14] = 0
$data[11,15] = 0
$data[11,16] = 0
$data[12,0] = 'Sixfold Demon Slash'
$data[12,1] = 10
$data[12,2] = 10
$data[12,3] = 9
$data[12,4] = 9
$data[12,5] = 16
$data[12,6] = 16
$data[12,7] = 0
$data[12,8] = 0
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0
$data[12,13] = 0
$data[12,14] = 0
$data[12,15] = 0
$data[12,16] = 0
$data[13,0] = 'Demon Flurry Rush'
$data[13,1] = 5
$data[13,2] = 5
$data[13,3] = 4
$data[13,4] = 4
$data[13,5] = 8
$data[13,6] = 8
$data[13,7] = 0
$data[13,8] = 0
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 0
$data[13,14] = 0
$data[13,15] = 0
$data[13,16] = 0
$data[14,0] = 'Right Fade Slash'
$data[14,1] = 7
$data[14,2] = 0
$data[14,3] = 0
$data[14,4] = 0
$data[14,5] = 0
$data[14,6] = 0
$data[14,7] = 0
$data[14,8] = 0
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 0
$data[14,13] = 0
$data[14,14] = 0
$data[14,15] = 0
$data[14,16] = 0
$data[15,0] = 'Left Fade Slash'
$data[15,1] = 7
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 0
$data[15,5] = 0
$data[15,6] = 0
$data[15,7] = 0
$data[15,8] = 0
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 0
$data[15,14] = 0
$data[15,15] = 0
$data[15,16] = 0
$data[16,0] = 'Right Double Round Slash'
$data[16,1] = 19
$data[16,2] = 7
$data[16,3] = 11
$data[16,4] = 0
$data[16,5] = 0
$data[16,6] = 0
$data[16,7] = 0
$data[16,8] = 0
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 0
$data[16,14] = 0
$data[16,15] = 0
$data[16,16] = 0
$data[17,0] = 'Left Double Round Slash'
$data[17,1] = 19
$data[17,2] = 7
$data[17,3] = 11
$data[17,4] = 0
$data[17,5] = 0
$data[17,6] = 0
$data[17,7] = 0
$data[17,8] = 0
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 0
$data[17,14] = 0
$data[17,15] = 0
$data[17,16] = 0
$data[18,0] = 'Blade Dance'
$data[18,1] = 17
$data[18,2] = 17
$data[18,3] = 6
$data[18,4] = 6
$data[18,5] = 10
$data[18,6] = 10
$data[18,7] = 9
$data[18,8] = 9
$data[18,9] = 11
$data[18,10] = 11
$data[18,11] = 9
$data[18,12] = 9
$data[18,13] = 12
$data[18,14] = 7
$data[18,15] = 12
$data[18,16] = 0
$data[19,0] = 'Demon Flurry'
$data[19,1] = 9
$data[19,2] = 9
$data[19,3] = 7
$data[19,4] = 7
$data[19,5] = 3
$data[19,6] = 17
$data[19,7] = 17
$data[19,8] = 0
$data[19,9] = 0
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 0
$data[19,13] = 0
$data[19,14] = 0
$data[19,15] = 0
$data[19,16] = 0
$data[20,0] = 'Jumping Double Slash'
$data[20,1] = 9
$data[20,2] = 12
$data[20,3] = 0
$data[20,4] = 0
$data[20,5] = 0
$data[20,6] = 0
$data[20,7] = 0
$data[20,8] = 0
$data[20,9] = 0
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 0
$data[20,14] = 0
$data[20,15] = 0
$data[20,16] = 0
$data[21,0] = 'Heavenly Blade Dance'
$data[21,1] = 17
$data[21,2] = 11
$data[21,3] = 11
$data[21,4] = 17
$data[21,5] = 20
$data[21,6] = 20
$data[21,7] = 12
$data[21,8] = 12
$data[21,9] = 21
$data[21,10] = 21
$data[21,11] = 0
$data[21,12] = 0
$data[21,13] = 0
$data[21,14] = 0
$data[21,15] = 0
$data[21,16] = 0
$data[22,0] = 'Midair Spinning Blade Dance'
$data[22,1] = 15
$data[22,2] = 10
$data[22,3] = 15
$data[22,4] = 10
$data[22,5] = 0
$data[22,6] = 0
$data[22,7] = 0
$data[22,8] = 0
$data[22,9] = 0
$data[22,10] = 0
$data[22,11] = 0
$data[22,12] = 0
$data[22,13] = 0
$data[22,14] = 0
$data[22,15] = 0
$data[22,16] = 0
$data[23,0] = 'Spinning Blade Dance Finisher'
$data[23,1] = 20
$data[23,2] = 20
$data[23,3] = 20
$data[23,4] = 20
$data[23,5] = 12
$data[23,6] = 12
$data[23,7] = 0
$data[23,8] = 0
$data[23,9] = 0
$data[23,10] = 0
$data[23,11] = 0
$data[23,12] = 0
$data[23,13] = 0
$data[23,14] = 0
$data[23,15] = 0
$data[23,16] = 0

$ws.Range("A1:Q24").Value = $data

# View: zoom + active selection as left by the editor
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("A22").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

